$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells store numeric-looking strings (prices / percentages) as
# plain TEXT, not numbers. Format each one as Text before writing so
# Excel does not auto-convert the new value into a number/percentage.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "297.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.60%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.75"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.79%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.011"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.36%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07531"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.76%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.608"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.49%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9177"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.20%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1181"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.29%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1825"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.83%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08954"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.94%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-7.21%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1051"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.39%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001278"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.88%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005948"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.32%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.08%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.376"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.19%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.22%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.263"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.71%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1372"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.01%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3222"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "17.48%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04092"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.83%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001266"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.35%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003927"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "7.35%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001301"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.59%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02405"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.96%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05202"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.41%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006305"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.36%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007837"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.19%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1324"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.82%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007415"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.68%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007090"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.47%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.78%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006588"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.93%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04543"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "27.04%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004205"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.08%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.03%"
